# Update for "Add data for 2022-03-02": refresh the "through" date from
# Feb 21 to Feb 22, 2022, and bump several neighborhood/month counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab
$ws.Name = "Through 2022-02-22"

# Update the header label for the "through" column (column B)
$ws.Range("B1").Value = "February 2022 (through February 22)"

# Row 3 - Austin
$ws.Range("D3").Value = 11
$ws.Range("F3").Value = 8
$ws.Range("H3").Value = 3

# Row 7 - Auburn Gresham
$ws.Range("D7").Value = 5
$ws.Range("L7").Value = 1

# Row 8 - North Lawndale
$ws.Range("N8").Value = 1

# Row 9 - United Center
$ws.Range("F9").Value = 1

# Row 23 - Washington Heights
$ws.Range("B23").Value = 3

# Row 32 - Grand Boulevard
$ws.Range("B32").Value = 1

# Row 36 - South Chicago
$ws.Range("L36").Value = 1

# Row 37 - West Pullman
$ws.Range("B37").Value = 2
$ws.Range("J37").Value = 1

# Row 38 - Wicker Park
$ws.Range("H38").Value = 1

# Row 39 - Rogers Park
$ws.Range("B39").Value = 5

# Row 41 - Loop
$ws.Range("B41").Value = 1

# Row 42 - River North
$ws.Range("L42").Value = 2

# Row 54 - Belmont Cragin
$ws.Range("J54").Value = 1

# Row 56 - Ashburn
$ws.Range("D56").Value = 1
$ws.Range("J56").Value = 2

# Row 67 - Garfield Ridge
$ws.Range("L67").Value = 1
